$d = $word.ActiveDocument

# Locate the paragraph that ends the list ("Implementar tabla de
# puntuación de las piezas.") and append a new list item after it,
# inheriting the same list style/formatting.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$newPara.Range.Text = "Añadir gestor de sonido y sonidos para el movimiento de las piezas y el ataque."
